$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the random-number formulas in B2:AK2 with plain sequential
# values 1..36 (no formula left behind).
for ($i = 1; $i -le 36; $i++) {
    $col = $i + 1
    $ws.Cells.Item(2, $col).Value = $i
}

# Move the active selection to AH19 (matches the saved view state).
$null = $ws.Range("AH19").Select()
